# Add a new "deals" worksheet after "contacts" with a couple of cells,
# then make it the active/selected sheet (mirrors the author's commit
# "Added New Test Cases for Deals Page").

$wb = $excel.ActiveWorkbook

# Worksheets.Add() drops the new sheet in front of the active sheet, so
# create it, fill it in, then move it to the end (after "contacts").
$deals = $wb.Worksheets.Add()
$deals.Name = "deals"

$deals.Range("A1").Value = "title"
$deals.Range("A2").Value = "abcd"

$deals.Move($null, $wb.Worksheets.Item("contacts"))

# Re-fetch the sheet by name before activating/selecting - the handle
# captured before the Move() call doesn't track tab-selection state.
$deals = $wb.Worksheets.Item("deals")
$deals.Activate()
$deals.Range("A2").Select()
